# Clarify the "Usage Rate" / "Lead Time" column headers in Example.xlsx.
#
# Renaming the header cells (rather than poking the ListObject/Table API
# directly) is what actually propagates through to:
#   - the shared-strings table (new text added, old unused text dropped)
#   - the Table1 column definitions (xl/tables/table1.xml), since the
#     table's header row is bound to these worksheet cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Usage Rate (Per Month)"
$ws.Range("B1").Value = "Lead Time (Days)"

# Widen columns A and B so the longer header text keeps fitting (matches
# the bestFit-style column autosize Excel performs after a header rename).
$ws.Columns.Item(1).ColumnWidth = 24.85546875
$ws.Columns.Item(2).ColumnWidth = 18.42578125
